$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H100").Value = 3070.389
$ws_ALC.Range("I100").Value = 1702
$ws_ALC.Range("K100").Value = 1702
$ws_ALC.Range("M100").Value = -1161
$ws_ALC.Range("H101").Value = 2628.4
$ws_ALC.Range("J101").Value = 4095
$ws_ALC.Range("L101").Value = 12285
$ws_ALC.Range("N101").Value = -15529
$ws_ALC.Range("H103").Value = 984
$ws_ALC.Range("I103").Value = 900.75
$ws_ALC.Range("K103").Value = 2702.25
$ws_ALC.Range("M103").Value = -2116.25
$ws_ALC.Range("H123").Value = 106666.664
$ws_ALC.Range("J123").Value = 106666.664
$ws_ALC.Range("L123").Value = 106666.664
$ws_ALC.Range("N123").Value = -116466.664
$ws_ALC.Range("H125").Value = 2793.1
$ws_ALC.Range("I125").Value = 1392.6
$ws_ALC.Range("J125").Value = 3259.9333
$ws_ALC.Range("K125").Value = 12533.4
$ws_ALC.Range("L125").Value = 29339.3997
$ws_ALC.Range("M125").Value = -10073.4
$ws_ALC.Range("N125").Value = -34259.3997
$ws_ALC.Range("H132").Value = 3265.5193
$ws_ALC.Range("I132").Value = 3319.745
$ws_ALC.Range("K132").Value = 9959.235000000001
$ws_ALC.Range("M132").Value = -7429.235000000001
$ws_ALC.Range("H137").Value = 26165.2
$ws_ALC.Range("I137").Value = 31061.586
$ws_ALC.Range("K137").Value = 93184.758
$ws_ALC.Range("M137").Value = -90634.758
$ws_ALC.Range("H138").Value = 2917.1516
$ws_ALC.Range("I138").Value = 1425.04
$ws_ALC.Range("J138").Value = 3421.2432
$ws_ALC.Range("K138").Value = 4275.12
$ws_ALC.Range("L138").Value = 10263.7296
$ws_ALC.Range("M138").Value = 864.8800000000001
$ws_ALC.Range("N138").Value = -20543.7296
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H5").Value = 401
$ws_ARM.Range("I5").Value = 434.77777
$ws_ARM.Range("K5").Value = 434.77777
$ws_ARM.Range("M5").Value = -322.77777
$ws_ARM.Range("H32").Value = 11064661
$ws_ARM.Range("I32").Value = 10837048
$ws_ARM.Range("K32").Value = 10837048
$ws_ARM.Range("M32").Value = -10836761
$ws_ARM.Range("H45").Value = 5599.636
$ws_ARM.Range("I45").Value = 5199.5
$ws_ARM.Range("K45").Value = 5199.5
$ws_ARM.Range("M45").Value = -4822.5
$ws_ARM.Range("H63").Value = 3847.5
$ws_ARM.Range("I63").Value = 2245
$ws_ARM.Range("K63").Value = 2245
$ws_ARM.Range("M63").Value = -1559
$ws_ARM.Range("H66").Value = 3847.5
$ws_ARM.Range("I66").Value = 2245
$ws_ARM.Range("K66").Value = 11225
$ws_ARM.Range("M66").Value = -7793
$ws_ARM.Range("H74").Value = 2669.8147
$ws_ARM.Range("I74").Value = 2547.2917
$ws_ARM.Range("K74").Value = 2547.2917
$ws_ARM.Range("M74").Value = -1673.2917
$ws_ARM.Range("H77").Value = 2669.8147
$ws_ARM.Range("I77").Value = 2547.2917
$ws_ARM.Range("K77").Value = 12736.4585
$ws_ARM.Range("M77").Value = -8368.458500000001
$ws_ARM.Range("H102").Value = 1918
$ws_ARM.Range("I102").Value = 1630.3636
$ws_ARM.Range("K102").Value = 1630.3636
$ws_ARM.Range("M102").Value = -8.363599999999906
$ws_ARM.Range("H109").Value = 63311.2
$ws_ARM.Range("I109").Value = 59999
$ws_ARM.Range("J109").Value = 64139.25
$ws_ARM.Range("K109").Value = 59999
$ws_ARM.Range("L109").Value = 64139.25
$ws_ARM.Range("M109").Value = -58612
$ws_ARM.Range("N109").Value = -66913.25
$ws_ARM.Range("H110").Value = 1902.2632
$ws_ARM.Range("I110").Value = 1649.4667
$ws_ARM.Range("J110").Value = 2850.25
$ws_ARM.Range("K110").Value = 1649.4667
$ws_ARM.Range("L110").Value = 2850.25
$ws_ARM.Range("M110").Value = 395.5333000000001
$ws_ARM.Range("N110").Value = -6940.25
$ws_ARM.Range("H121").Value = 37830.6
$ws_ARM.Range("J121").Value = 37830.6
$ws_ARM.Range("L121").Value = 37830.6
$ws_ARM.Range("N121").Value = -41324.6
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H4").Value = 401
$ws_BSM.Range("I4").Value = 434.77777
$ws_BSM.Range("K4").Value = 434.77777
$ws_BSM.Range("M4").Value = -319.77777
$ws_BSM.Range("H107").Value = 1693.25
$ws_BSM.Range("I107").Value = 924.3889
$ws_BSM.Range("J107").Value = 3999.8333
$ws_BSM.Range("K107").Value = 924.3889
$ws_BSM.Range("L107").Value = 3999.8333
$ws_BSM.Range("M107").Value = 995.6111
$ws_BSM.Range("N107").Value = -7839.8333
$ws_BSM.Range("H108").Value = 99989.5
$ws_BSM.Range("J108").Value = 99989.5
$ws_BSM.Range("L108").Value = 99989.5
$ws_BSM.Range("N108").Value = -107669.5
$ws_BSM.Range("H123").Value = 94992
$ws_BSM.Range("J123").Value = 94992
$ws_BSM.Range("L123").Value = 94992
$ws_BSM.Range("N123").Value = -104792
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 3736.9167
$ws_CRP.Range("I31").Value = 2742.25
$ws_CRP.Range("J31").Value = 4980.25
$ws_CRP.Range("K31").Value = 2742.25
$ws_CRP.Range("L31").Value = 4980.25
$ws_CRP.Range("M31").Value = -2447.25
$ws_CRP.Range("N31").Value = -5570.25
$ws_CRP.Range("H34").Value = 3736.9167
$ws_CRP.Range("I34").Value = 2742.25
$ws_CRP.Range("J34").Value = 4980.25
$ws_CRP.Range("K34").Value = 2742.25
$ws_CRP.Range("L34").Value = 4980.25
$ws_CRP.Range("M34").Value = -2540.25
$ws_CRP.Range("N34").Value = -5384.25
$ws_CRP.Range("H58").Value = 2669.4773
$ws_CRP.Range("I58").Value = 2432.3713
$ws_CRP.Range("J58").Value = 3591.5557
$ws_CRP.Range("K58").Value = 2432.3713
$ws_CRP.Range("L58").Value = 3591.5557
$ws_CRP.Range("M58").Value = -2229.3713
$ws_CRP.Range("N58").Value = -3997.5557
$ws_CRP.Range("H123").Value = 62570
$ws_CRP.Range("J123").Value = 62570
$ws_CRP.Range("L123").Value = 62570
$ws_CRP.Range("N123").Value = -72370
$ws_CRP.Range("H136").Value = 2669.4773
$ws_CRP.Range("I136").Value = 2432.3713
$ws_CRP.Range("J136").Value = 3591.5557
$ws_CRP.Range("K136").Value = 7297.113899999999
$ws_CRP.Range("L136").Value = 10774.6671
$ws_CRP.Range("M136").Value = -4747.113899999999
$ws_CRP.Range("N136").Value = -15874.6671
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 65299852
$ws_CUL.Range("I4").Value = 58400450
$ws_CUL.Range("K4").Value = 175201350
$ws_CUL.Range("M4").Value = -175201238
$ws_CUL.Range("H56").Value = 9699.799999999999
$ws_CUL.Range("I56").Value = 9699.799999999999
$ws_CUL.Range("K56").Value = 9699.799999999999
$ws_CUL.Range("M56").Value = -9169.799999999999
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 2102.5
$ws_GSM.Range("I80").Value = 905
$ws_GSM.Range("J80").Value = 2342
$ws_GSM.Range("K80").Value = 905
$ws_GSM.Range("L80").Value = 2342
$ws_GSM.Range("M80").Value = 93
$ws_GSM.Range("N80").Value = -4338
$ws_GSM.Range("H83").Value = 2102.5
$ws_GSM.Range("I83").Value = 905
$ws_GSM.Range("J83").Value = 2342
$ws_GSM.Range("K83").Value = 4525
$ws_GSM.Range("L83").Value = 11710
$ws_GSM.Range("M83").Value = 467
$ws_GSM.Range("N83").Value = -21694
$ws_GSM.Range("H102").Value = 1497.5
$ws_GSM.Range("I102").Value = 1497.5
$ws_GSM.Range("J102").Value = 0
$ws_GSM.Range("K102").Value = 1497.5
$ws_GSM.Range("L102").Value = 0
$ws_GSM.Range("M102").Value = 124.5
$ws_GSM.Range("N102").ClearContents()
$ws_GSM.Range("H107").Value = 634.087
$ws_GSM.Range("I107").Value = 649.3333
$ws_GSM.Range("K107").Value = 649.3333
$ws_GSM.Range("M107").Value = 1270.6667
$ws_GSM.Range("H122").Value = 1762.7142
$ws_GSM.Range("I122").Value = 1556.5
$ws_GSM.Range("J122").Value = 3000
$ws_GSM.Range("K122").Value = 4669.5
$ws_GSM.Range("L122").Value = 9000
$ws_GSM.Range("M122").Value = -2219.5
$ws_GSM.Range("N122").Value = -13900
$ws_GSM.Range("H124").Value = 81332.336
$ws_GSM.Range("J124").Value = 81332.336
$ws_GSM.Range("L124").Value = 81332.336
$ws_GSM.Range("N124").Value = -91152.336
$ws_GSM.Range("H130").Value = 100000
$ws_GSM.Range("J130").Value = 100000
$ws_GSM.Range("L130").Value = 100000
$ws_GSM.Range("N130").Value = -110040
$ws_GSM.Range("H132").Value = 4470.4287
$ws_GSM.Range("I132").Value = 4493.684
$ws_GSM.Range("K132").Value = 13481.052
$ws_GSM.Range("M132").Value = -10951.052
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H61").Value = 5296.067
$ws_LTW.Range("I61").Value = 2854
$ws_LTW.Range("J61").Value = 15064.333
$ws_LTW.Range("K61").Value = 2854
$ws_LTW.Range("L61").Value = 15064.333
$ws_LTW.Range("M61").Value = -2652
$ws_LTW.Range("N61").Value = -15468.333
$ws_LTW.Range("H74").Value = 34881.824
$ws_LTW.Range("I74").Value = 33995.5
$ws_LTW.Range("K74").Value = 33995.5
$ws_LTW.Range("M74").Value = -32997.5
$ws_LTW.Range("H77").Value = 34881.824
$ws_LTW.Range("I77").Value = 33995.5
$ws_LTW.Range("K77").Value = 101986.5
$ws_LTW.Range("M77").Value = -96994.5
$ws_LTW.Range("H93").Value = 30304430
$ws_LTW.Range("I93").Value = 66667696
$ws_LTW.Range("J93").Value = 1709.3889
$ws_LTW.Range("K93").Value = 66667696
$ws_LTW.Range("L93").Value = 1709.3889
$ws_LTW.Range("M93").Value = -66666448
$ws_LTW.Range("N93").Value = -4205.3889
$ws_LTW.Range("H113").Value = 5296.067
$ws_LTW.Range("I113").Value = 2854
$ws_LTW.Range("J113").Value = 15064.333
$ws_LTW.Range("K113").Value = 2854
$ws_LTW.Range("L113").Value = 15064.333
$ws_LTW.Range("M113").Value = -684
$ws_LTW.Range("N113").Value = -19404.333
$ws_LTW.Range("H132").Value = 318154.4
$ws_LTW.Range("I132").Value = 345259.1
$ws_LTW.Range("K132").Value = 1035777.3
$ws_LTW.Range("M132").Value = -1033247.3
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H81").Value = 47120.652
$ws_WVR.Range("I81").Value = 79873.62
$ws_WVR.Range("J81").Value = 4541.8
$ws_WVR.Range("K81").Value = 159747.24
$ws_WVR.Range("L81").Value = 9083.6
$ws_WVR.Range("M81").Value = -158686.24
$ws_WVR.Range("N81").Value = -11205.6
$ws_WVR.Range("H84").Value = 47120.652
$ws_WVR.Range("I84").Value = 79873.62
$ws_WVR.Range("J84").Value = 4541.8
$ws_WVR.Range("K84").Value = 798736.2
$ws_WVR.Range("L84").Value = 45418
$ws_WVR.Range("M84").Value = -793432.2
$ws_WVR.Range("N84").Value = -56026
$ws_WVR.Range("H107").Value = 359
$ws_WVR.Range("I107").Value = 380.6
$ws_WVR.Range("K107").Value = 1141.8
$ws_WVR.Range("M107").Value = 778.1999999999998
$ws_WVR.Range("H124").Value = 31000
$ws_WVR.Range("J124").Value = 31000
$ws_WVR.Range("L124").Value = 31000
$ws_WVR.Range("N124").Value = -40820
$ws_WVR.Range("H131").Value = 107398.2
$ws_WVR.Range("J131").Value = 107398.2
$ws_WVR.Range("L131").Value = 107398.2
$ws_WVR.Range("N131").Value = -117478.2
$ws_WVR.Range("H136").Value = 1652.0317
$ws_WVR.Range("I136").Value = 1297.4897
$ws_WVR.Range("J136").Value = 2892.9285
$ws_WVR.Range("K136").Value = 3892.4691
$ws_WVR.Range("L136").Value = 8678.7855
$ws_WVR.Range("M136").Value = -1342.4691
$ws_WVR.Range("N136").Value = -13778.7855
